$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Range("G12").Value = "共用代碼檔`n0:未處理`n1:兌現入帳`n2:退票`n3:抽票`n4:兌現未入帳`n5:即期票"
